# GUI spectra merging started + multiple bug fixes
$wb = $excel.ActiveWorkbook

# Rename "constant_names" sheet to "target"
$wsTarget = $wb.Worksheets.Item("constant_names")
$wsTarget.Name = "target"

# Populate the renamed "target" sheet with merged spectra constant/wavelength info
$wsTarget.Range("A1").Value = "constant"
$wsTarget.Range("B1").Value = "SB"
$wsTarget.Range("A2").Value = "wavelength"
$wsTarget.Range("B2").Value = 306
$wsTarget.Range("C2").Value = 387

# Make "target" sheet the active/selected sheet, with B1 selected
$wsTarget.Select()
$wsTarget.Range("B1").Select()
